$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 84.11458333333334
$ws.Range("C2").Value = 84.11458333333334
$ws.Range("D2").Value = 84.11458333333334
$ws.Range("E2:U2").Value = 83.85416666666666
